# JS-SPA-Self-Evaluation-Protocol.xlsx
# "update edit ad and edit profile"
#
# Fills in the self-evaluation scores for the "Edit User Profile" (row 29),
# "Change Password" (row 30) and "Logout" (row 31) checklist items, and
# moves the active selection/scroll position to reflect where the user was
# last working in the sheet (around E31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Score the remaining "Basic Options" checklist rows.
$ws.Range("E29").Value = 5
$ws.Range("E30").Value = 5
$ws.Range("E31").Value = 3

# Move the view/selection to where the edits were made.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("E31").Select()
